$wb = $excel.ActiveWorkbook

# Add the new worksheet "VATRIM" after the last existing sheet
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "VATRIM"

# Match the page setup used by the other sheets in the workbook
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Header row (bold, orange-filled style), matching the style used on the other sheets
$ws.Range("A1").Value = "time"
$ws.Range("B1").Value = "VA_trim"
$headerSrc = $wb.Worksheets.Item("PREVISION").Range("A1")
$headerSrc.Copy()
$ws.Range("A1:B1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data rows 2..101 (quarterly "time" / "VA_trim" series)
$data = New-Object "object[,]" 100,2
$data[0,0] = "'1997-01-01"
$data[0,1] = 132947.420754456
$data[1,0] = "'1997-04-01"
$data[1,1] = 212179.739504027
$data[2,0] = "'1997-07-01"
$data[2,1] = 212276.33848145
$data[3,0] = "'1997-10-01"
$data[3,1] = 227150.506767315
$data[4,0] = "'1998-01-01"
$data[4,1] = 141393.103621597
$data[5,0] = "'1998-04-01"
$data[5,1] = 226795.429710672
$data[6,0] = "'1998-07-01"
$data[6,1] = 226141.013334937
$data[7,0] = "'1998-10-01"
$data[7,1] = 239206.415486733
$data[8,0] = "'1999-01-01"
$data[8,1] = 151583.382586458
$data[9,0] = "'1999-04-01"
$data[9,1] = 241196.029205423
$data[10,0] = "'1999-07-01"
$data[10,1] = 240580.767817142
$data[11,0] = "'1999-10-01"
$data[11,1] = 257163.48842593
$data[12,0] = "'2000-01-01"
$data[12,1] = 162331.055571178
$data[13,0] = "'2000-04-01"
$data[13,1] = 255636.177730075
$data[14,0] = "'2000-07-01"
$data[14,1] = 254761.346080887
$data[15,0] = "'2000-10-01"
$data[15,1] = 268794.802071233
$data[16,0] = "'2001-01-01"
$data[16,1] = 172152.944149418
$data[17,0] = "'2001-04-01"
$data[17,1] = 262959.217922613
$data[18,0] = "'2001-07-01"
$data[18,1] = 259489.313219138
$data[19,0] = "'2001-10-01"
$data[19,1] = 272553.853557768
$data[20,0] = "'2002-01-01"
$data[20,1] = 165954.121361025
$data[21,0] = "'2002-04-01"
$data[21,1] = 254451.023334684
$data[22,0] = "'2002-07-01"
$data[22,1] = 251850.110041742
$data[23,0] = "'2002-10-01"
$data[23,1] = 265605.637329344
$data[24,0] = "'2003-01-01"
$data[24,1] = 165394.611610315
$data[25,0] = "'2003-04-01"
$data[25,1] = 260338.040947399
$data[26,0] = "'2003-07-01"
$data[26,1] = 259586.928674046
$data[27,0] = "'2003-10-01"
$data[27,1] = 276279.52476645
$data[28,0] = "'2004-01-01"
$data[28,1] = 169633.912671073
$data[29,0] = "'2004-04-01"
$data[29,1] = 266450.157233896
$data[30,0] = "'2004-07-01"
$data[30,1] = 266503.51816605
$data[31,0] = "'2004-10-01"
$data[31,1] = 283239.162626905
$data[32,0] = "'2005-01-01"
$data[32,1] = 181107.650782012
$data[33,0] = "'2005-04-01"
$data[33,1] = 276683.096726239
$data[34,0] = "'2005-07-01"
$data[34,1] = 274530.816452221
$data[35,0] = "'2005-10-01"
$data[35,1] = 288409.591236635
$data[36,0] = "'2006-01-01"
$data[36,1] = 181453.27855698
$data[37,0] = "'2006-04-01"
$data[37,1] = 270137.047097937
$data[38,0] = "'2006-07-01"
$data[38,1] = 263887.320242544
$data[39,0] = "'2006-10-01"
$data[39,1] = 272851.771650393
$data[40,0] = "'2007-01-01"
$data[40,1] = 169545.058154116
$data[41,0] = "'2007-04-01"
$data[41,1] = 270276.397914076
$data[42,0] = "'2007-07-01"
$data[42,1] = 268450.180533682
$data[43,0] = "'2007-10-01"
$data[43,1] = 288675.05677442
$data[44,0] = "'2008-01-01"
$data[44,1] = 171154.950219043
$data[45,0] = "'2008-04-01"
$data[45,1] = 288532.41522993
$data[46,0] = "'2008-07-01"
$data[46,1] = 289609.170758743
$data[47,0] = "'2008-10-01"
$data[47,1] = 308698.106163459
$data[48,0] = "'2009-01-01"
$data[48,1] = 183405.31431924
$data[49,0] = "'2009-04-01"
$data[49,1] = 296357.683513839
$data[50,0] = "'2009-07-01"
$data[50,1] = 291988.326762593
$data[51,0] = "'2009-10-01"
$data[51,1] = 307280.861394363
$data[52,0] = "'2010-01-01"
$data[52,1] = 168979.639599388
$data[53,0] = "'2010-04-01"
$data[53,1] = 274326.632053725
$data[54,0] = "'2010-07-01"
$data[54,1] = 270052.111531832
$data[55,0] = "'2010-10-01"
$data[55,1] = 285782.401860167
$data[56,0] = "'2011-01-01"
$data[56,1] = 153859.834863015
$data[57,0] = "'2011-04-01"
$data[57,1] = 262279.881650588
$data[58,0] = "'2011-07-01"
$data[58,1] = 262568.888808087
$data[59,0] = "'2011-10-01"
$data[59,1] = 283600.019780695
$data[60,0] = "'2012-01-01"
$data[60,1] = 158488.861968795
$data[61,0] = "'2012-04-01"
$data[61,1] = 264177.013235019
$data[62,0] = "'2012-07-01"
$data[62,1] = 257066.434451383
$data[63,0] = "'2012-10-01"
$data[63,1] = 263471.781077842
$data[64,0] = "'2013-01-01"
$data[64,1] = 147564.46110679
$data[65,0] = "'2013-04-01"
$data[65,1] = 270161.370594591
$data[66,0] = "'2013-07-01"
$data[66,1] = 270234.255710547
$data[67,0] = "'2013-10-01"
$data[67,1] = 302423.029869192
$data[68,0] = "'2014-01-01"
$data[68,1] = 174444.054917061
$data[69,0] = "'2014-04-01"
$data[69,1] = 338611.306219553
$data[70,0] = "'2014-07-01"
$data[70,1] = 342392.336499563
$data[71,0] = "'2014-10-01"
$data[71,1] = 365941.859530217
$data[72,0] = "'2015-01-01"
$data[72,1] = 203613.060426295
$data[73,0] = "'2015-04-01"
$data[73,1] = 348369.318706765
$data[74,0] = "'2015-07-01"
$data[74,1] = 341264.200852872
$data[75,0] = "'2015-10-01"
$data[75,1] = 362592.580291987
$data[76,0] = "'2016-01-01"
$data[76,1] = 197030.448484361
$data[77,0] = "'2016-04-01"
$data[77,1] = 347853.535772949
$data[78,0] = "'2016-07-01"
$data[78,1] = 345871.87410849
$data[79,0] = "'2016-10-01"
$data[79,1] = 371745.201973878
$data[80,0] = "'2017-01-01"
$data[80,1] = 211599.878486715
$data[81,0] = "'2017-04-01"
$data[81,1] = 374923.023234774
$data[82,0] = "'2017-07-01"
$data[82,1] = 374356.279947385
$data[83,0] = "'2017-10-01"
$data[83,1] = 402250.239172714
$data[84,0] = "'2018-01-01"
$data[84,1] = 218161.644542643
$data[85,0] = "'2018-04-01"
$data[85,1] = 374051.031218892
$data[86,0] = "'2018-07-01"
$data[86,1] = 370119.002327364
$data[87,0] = "'2018-10-01"
$data[87,1] = 392365.714838474
$data[88,0] = "'2019-01-01"
$data[88,1] = 224956.232518354
$data[89,0] = "'2019-04-01"
$data[89,1] = 383870.997190286
$data[90,0] = "'2019-07-01"
$data[90,1] = 379616.057817982
$data[91,0] = "'2019-10-01"
$data[91,1] = 405535.752239579
$data[92,0] = "'2020-01-01"
$data[92,1] = 225694.858612034
$data[93,0] = "'2020-04-01"
$data[93,1] = 391619.225768974
$data[94,0] = "'2020-07-01"
$data[94,1] = 389280.648419773
$data[95,0] = "'2020-10-01"
$data[95,1] = 415879.343420858
$data[96,0] = "'2021-01-01"
$data[96,1] = 232721.751457131
$data[97,0] = "'2021-04-01"
$data[97,1] = 403904.474216792
$data[98,0] = "'2021-07-01"
$data[98,1] = 401534.97246494
$data[99,0] = "'2021-10-01"
$data[99,1] = 429647.176254566
$ws.Range("A2:B101").Value = $data
$ws.Range("A2:A101").Style = "Normal"
